# Reorders data rows 2-11 on the active sheet according to the mapping
# observed between the original and the updated workbook. Every column
# (A:AY) for a given source row travels together as a unit - i.e. whole
# records were shuffled, not individual fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> original (source) row number
$mapping = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 10
    6  = 4
    7  = 11
    8  = 9
    9  = 3
    10 = 8
    11 = 7
}

# Snapshot every source row (columns A:AY) before writing anything back,
# since several destinations read from rows that are also overwritten.
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("A$($srcRow):AY$($srcRow)").Value2
    }
}

foreach ($newRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$newRow]
    # Columns Y (Startdatum) and AA (Slutdatum) hold plain text that looks
    # like a date ("2023-07-03"). Force text formatting first so the bulk
    # Value2 array write below does not auto-coerce those cells into real
    # date serial numbers.
    $ws.Range("Y$($newRow):AA$($newRow)").NumberFormat = "@"
    $ws.Range("A$($newRow):AY$($newRow)").Value2 = $snapshot[$srcRow]
}
